$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.03803729543516
$ws.Range("D2").Value = 1.040681661590539
$ws.Range("E2").Value = 1.036634424138333
$ws.Range("F2").Value = 1.036750166034662
$ws.Range("I2").Value = 1.033126696701326
$ws.Range("J2").Value = 1.043136909830661
$ws.Range("K2").Value = 1.043463319468093
$ws.Range("L2").Value = 1.039427604399025
$ws.Range("M2").Value = 1.039543015424148
$ws.Range("N2").Value = 1.044618283892154
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.039948823371063
$ws.Range("D3").Value = 1.04250774979497
$ws.Range("E3").Value = 1.038292428549118
$ws.Range("F3").Value = 1.039302666831624
$ws.Range("I3").Value = 1.033520507920174
$ws.Range("J3").Value = 1.04468824180698
$ws.Range("K3").Value = 1.045097212600779
$ws.Range("L3").Value = 1.040893002758833
$ws.Range("M3").Value = 1.041900569223949
$ws.Range("N3").Value = 1.046171818937762
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.04118032987849
$ws.Range("D4").Value = 1.043684353180134
$ws.Range("E4").Value = 1.039360631061186
$ws.Range("F4").Value = 1.040948089337549
$ws.Range("I4").Value = 1.03377193629464
$ws.Range("J4").Value = 1.045686511083079
$ws.Range("K4").Value = 1.046149039009968
$ws.Range("L4").Value = 1.041836126670144
$ws.Range("M4").Value = 1.043419604892013
$ws.Range("N4").Value = 1.047171505870737
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.041696794668497
$ws.Range("D5").Value = 1.044177826151174
$ws.Range("E5").Value = 1.039808615855635
$ws.Range("F5").Value = 1.041638378554419
$ws.Range("I5").Value = 1.033876831327163
$ws.Range("J5").Value = 1.046104878785228
$ws.Range("K5").Value = 1.046589955069662
$ws.Range("L5").Value = 1.042231419455727
$ws.Range("M5").Value = 1.044056700030551
$ws.Range("N5").Value = 1.04759046770301
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.041783438251973
$ws.Range("D6").Value = 1.044260614494265
$ws.Range("E6").Value = 1.039883771353047
$ws.Range("F6").Value = 1.041754197539442
$ws.Range("I6").Value = 1.033894396635238
$ws.Range("J6").Value = 1.046175048721507
$ws.Range("K6").Value = 1.046663912927014
$ws.Range("L6").Value = 1.042297721304791
$ws.Range("M6").Value = 1.044163583833217
$ws.Range("N6").Value = 1.047660737288647
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.041187235822499
$ws.Range("D7").Value = 1.043690951557537
$ws.Range("E7").Value = 1.039366621298763
$ws.Range("F7").Value = 1.040957318643299
$ws.Range("I7").Value = 1.033773341062582
$ws.Range("J7").Value = 1.045692106428305
$ws.Range("K7").Value = 1.046154935517043
$ws.Range("L7").Value = 1.041841413264475
$ws.Range("M7").Value = 1.043428123659575
$ws.Range("N7").Value = 1.047177109161994
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.038684435803203
$ws.Range("D8").Value = 1.041299846861556
$ws.Range("E8").Value = 1.037195728363425
$ws.Range("F8").Value = 1.037614108045104
$ws.Range("I8").Value = 1.033260494023322
$ws.Range("J8").Value = 1.04366235211291
$ws.Range("K8").Value = 1.044016636262473
$ws.Range("L8").Value = 1.039923908797394
$ws.Range("M8").Value = 1.04034112278925
$ws.Range("N8").Value = 1.045144472362655
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.034231736301966
$ws.Range("D9").Value = 1.037046986628221
$ws.Range("E9").Value = 1.033333796630768
$ws.Range("F9").Value = 1.031673441845578
$ws.Range("I9").Value = 1.03233047994801
$ws.Range("J9").Value = 1.040042141170344
$ws.Range("K9").Value = 1.040206152779362
$ws.Range("L9").Value = 1.036505096507577
$ws.Range("M9").Value = 1.034850198310076
$ws.Range("N9").Value = 1.041519120305307
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.03123297791111
$ws.Range("D10").Value = 1.034183612054311
$ws.Range("E10").Value = 1.030733159331769
$ws.Range("F10").Value = 1.02767704409496
$ws.Range("I10").Value = 1.031692327480197
$ws.Range("J10").Value = 1.03759792537287
$ws.Range("K10").Value = 1.037635712946422
$ws.Range("L10").Value = 1.0341976818563
$ws.Range("M10").Value = 1.031152645573755
$ws.Range("N10").Value = 1.039071433441046
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.029926902592536
$ws.Range("D11").Value = 1.032936700026045
$ws.Range("E11").Value = 1.029600557674375
$ws.Range("F11").Value = 1.025937430667396
$ws.Range("I11").Value = 1.031411594804677
$ws.Range("J11").Value = 1.036531928158311
$ws.Range("K11").Value = 1.036515197017374
$ws.Range("L11").Value = 1.033191548289875
$ws.Range("M11").Value = 1.029542243132827
$ws.Range("N11").Value = 1.03800392238817
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.029440592951192
$ws.Range("D12").Value = 1.032472450070467
$ws.Range("E12").Value = 1.029178852770513
$ws.Range("F12").Value = 1.025289834366208
$ws.Range("I12").Value = 1.031306646918339
$ws.Range("J12").Value = 1.03613479316189
$ws.Range("K12").Value = 1.036097831046899
$ws.Range("L12").Value = 1.032816745895571
$ws.Range("M12").Value = 1.028942616404132
$ws.Range("N12").Value = 1.037606223414501
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.029544961754052
$ws.Range("D13").Value = 1.032572083162974
$ws.Range("E13").Value = 1.029269355898446
$ws.Range("F13").Value = 1.025428811367775
$ws.Range("I13").Value = 1.031329189118125
$ws.Range("J13").Value = 1.036220033661799
$ws.Range("K13").Value = 1.036187410282339
$ws.Range("L13").Value = 1.032897191569604
$ws.Range("M13").Value = 1.029071304849781
$ws.Range("N13").Value = 1.037691584965698
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.029886728258668
$ws.Range("D14").Value = 1.03289834742849
$ws.Range("E14").Value = 1.029565720108066
$ws.Range("F14").Value = 1.025883929559118
$ws.Range("I14").Value = 1.031402933525775
$ws.Range("J14").Value = 1.036499125018719
$ws.Range("K14").Value = 1.036480721182446
$ws.Range("L14").Value = 1.033160589186813
$ws.Range("M14").Value = 1.029492707696552
$ws.Range("N14").Value = 1.037971072664358
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.030097145130206
$ws.Range("D15").Value = 1.033099224039295
$ws.Range("E15").Value = 1.02974818568719
$ws.Range("F15").Value = 1.026164152304924
$ws.Range("I15").Value = 1.031448280668561
$ws.Range("J15").Value = 1.036670925684783
$ws.Range("K15").Value = 1.036661285562233
$ws.Range("L15").Value = 1.033322733287424
$ws.Range("M15").Value = 1.029752154081309
$ws.Range("N15").Value = 1.038143117307072
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.031319494933161
$ws.Range("D16").Value = 1.034266214218273
$ws.Range("E16").Value = 1.030808186834825
$ws.Range("F16").Value = 1.027792299360945
$ws.Range("I16").Value = 1.031710865142146
$ws.Range("J16").Value = 1.037668508759645
$ws.Range("K16").Value = 1.037709917373516
$ws.Range("L16").Value = 1.034264305709114
$ws.Range("M16").Value = 1.03125932181087
$ws.Range("N16").Value = 1.039142117064326
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.032084185723482
$ws.Range("D17").Value = 1.034996326023772
$ws.Range("E17").Value = 1.031471335147675
$ws.Range("F17").Value = 1.028811107568865
$ws.Range("I17").Value = 1.031874390911894
$ws.Range("J17").Value = 1.038292201679246
$ws.Range("K17").Value = 1.038365667900064
$ws.Range("L17").Value = 1.03485303429453
$ws.Range("M17").Value = 1.032202194256031
$ws.Range("N17").Value = 1.03976669569941
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.032529486644412
$ws.Range("D18").Value = 1.035421509655738
$ws.Range("E18").Value = 1.0318575123012
$ws.Range("F18").Value = 1.029404482196602
$ws.Range("I18").Value = 1.03196934798895
$ws.Range("J18").Value = 1.038655256462901
$ws.Range("K18").Value = 1.038747434635467
$ws.Range("L18").Value = 1.035195755397877
$ws.Range("M18").Value = 1.032751258179271
$ws.Range("N18").Value = 1.040130266062502
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.032681199912693
$ws.Range("D19").Value = 1.035566372201628
$ws.Range("E19").Value = 1.031989083395387
$ws.Range("F19").Value = 1.029606660089183
$ws.Range("I19").Value = 1.032001654123741
$ws.Range("J19").Value = 1.038778925006349
$ws.Range("K19").Value = 1.038877485752409
$ws.Range("L19").Value = 1.035312500863317
$ws.Range("M19").Value = 1.03293832399577
$ws.Range("N19").Value = 1.040254110229466
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.032002217421621
$ws.Range("D20").Value = 1.034918062318243
$ws.Range("E20").Value = 1.031400250567719
$ws.Range("F20").Value = 1.028701890380341
$ws.Range("I20").Value = 1.031856890141346
$ws.Range("J20").Value = 1.038225361537993
$ws.Range("K20").Value = 1.038295386912581
$ws.Range("L20").Value = 1.034789939229309
$ws.Range("M20").Value = 1.032101126160789
$ws.Range("N20").Value = 1.03969976063749
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.029786119258999
$ws.Range("D21").Value = 1.032802301118166
$ws.Range("E21").Value = 1.029478476242039
$ws.Range("F21").Value = 1.02574994839542
$ws.Range("I21").Value = 1.031381236238843
$ws.Range("J21").Value = 1.036416972273709
$ws.Range("K21").Value = 1.036394380599735
$ws.Range("L21").Value = 1.033083055172692
$ws.Range("M21").Value = 1.029368655526789
$ws.Range("N21").Value = 1.037888803253026
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.028385954085297
$ws.Range("D22").Value = 1.031465707914179
$ws.Range("E22").Value = 1.028264344345441
$ws.Range("F22").Value = 1.023885665932058
$ws.Range("I22").Value = 1.031078285427345
$ws.Range("J22").Value = 1.035273146481948
$ws.Range("K22").Value = 1.035192436339383
$ws.Range("L22").Value = 1.032003609648211
$ws.Range("M22").Value = 1.02764221975407
$ws.Range("N22").Value = 1.036743353097443
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.029128865863485
$ws.Range("D23").Value = 1.032174872064017
$ws.Range("E23").Value = 1.02890854143985
$ws.Range("F23").Value = 1.02487475981064
$ws.Range("I23").Value = 1.031239256985507
$ws.Range("J23").Value = 1.035880166421489
$ws.Range("K23").Value = 1.035830255519411
$ws.Range("L23").Value = 1.032576446554198
$ws.Range("M23").Value = 1.02855825091298
$ws.Range("N23").Value = 1.037351235074923
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.032039257624187
$ws.Range("D24").Value = 1.034953428412488
$ws.Range("E24").Value = 1.031432372563014
$ws.Range("F24").Value = 1.028751243640343
$ws.Range("I24").Value = 1.031864799298352
$ws.Range("J24").Value = 1.038255565989966
$ws.Range("K24").Value = 1.038327146093757
$ws.Range("L24").Value = 1.034818451256022
$ws.Range("M24").Value = 1.032146797254132
$ws.Range("N24").Value = 1.03973000798325
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.035388074676226
$ws.Range("D25").Value = 1.038151291011355
$ws.Range("E25").Value = 1.034336676422887
$ws.Range("F25").Value = 1.033215384898
$ws.Range("I25").Value = 1.032574075052542
$ws.Range("J25").Value = 1.040983361002991
$ws.Range("K25").Value = 1.041196449144188
$ws.Range("L25").Value = 1.037393811629086
$ws.Range("M25").Value = 1.036276059107942
$ws.Range("N25").Value = 1.042461676778076
